$d = $word.ActiveDocument

$d.Content.Find.Execute("53×18=954", $true, $false, $false, $false, $false, $true, 1, $false, "19×27=513", 2) | Out-Null
$d.Content.Find.Execute("45×25=1125", $true, $false, $false, $false, $false, $true, 1, $false, "51×87=4437", 2) | Out-Null
$d.Content.Find.Execute("76×18=1368", $true, $false, $false, $false, $false, $true, 1, $false, "97×34=3298", 2) | Out-Null
$d.Content.Find.Execute("72×77=5544", $true, $false, $false, $false, $false, $true, 1, $false, "62×82=5084", 2) | Out-Null
$d.Content.Find.Execute("43×72=3096", $true, $false, $false, $false, $false, $true, 1, $false, "75×87=6525", 2) | Out-Null
$d.Content.Find.Execute("75×84=6300", $true, $false, $false, $false, $false, $true, 1, $false, "38×35=1330", 2) | Out-Null
$d.Content.Find.Execute("38×61=2318", $true, $false, $false, $false, $false, $true, 1, $false, "82×67=5494", 2) | Out-Null
$d.Content.Find.Execute("70×47=3290", $true, $false, $false, $false, $false, $true, 1, $false, "85×67=5695", 2) | Out-Null
$d.Content.Find.Execute("56×61=3416", $true, $false, $false, $false, $false, $true, 1, $false, "62×60=3720", 2) | Out-Null
$d.Content.Find.Execute("78×74=5772", $true, $false, $false, $false, $false, $true, 1, $false, "66×87=5742", 2) | Out-Null
$d.Content.Find.Execute("61×83=5063", $true, $false, $false, $false, $false, $true, 1, $false, "84×30=2520", 2) | Out-Null
$d.Content.Find.Execute("72×63=4536", $true, $false, $false, $false, $false, $true, 1, $false, "60×56=3360", 2) | Out-Null
$d.Content.Find.Execute("72×22=1584", $true, $false, $false, $false, $false, $true, 1, $false, "95×82=7790", 2) | Out-Null
$d.Content.Find.Execute("31×46=1426", $true, $false, $false, $false, $false, $true, 1, $false, "63×71=4473", 2) | Out-Null
$d.Content.Find.Execute("54×46=2484", $true, $false, $false, $false, $false, $true, 1, $false, "71×78=5538", 2) | Out-Null
$d.Content.Find.Execute("35×20=700", $true, $false, $false, $false, $false, $true, 1, $false, "81×28=2268", 2) | Out-Null
$d.Content.Find.Execute("38×56=2128", $true, $false, $false, $false, $false, $true, 1, $false, "65×30=1950", 2) | Out-Null
$d.Content.Find.Execute("91×21=1911", $true, $false, $false, $false, $false, $true, 1, $false, "44×34=1496", 2) | Out-Null
$d.Content.Find.Execute("72×17=1224", $true, $false, $false, $false, $false, $true, 1, $false, "60×41=2460", 2) | Out-Null
$d.Content.Find.Execute("99×70=6930", $true, $false, $false, $false, $false, $true, 1, $false, "61×71=4331", 2) | Out-Null
$d.Content.Find.Execute("49×16=784", $true, $false, $false, $false, $false, $true, 1, $false, "66×81=5346", 2) | Out-Null
$d.Content.Find.Execute("42×42=1764", $true, $false, $false, $false, $false, $true, 1, $false, "76×41=3116", 2) | Out-Null
$d.Content.Find.Execute("28×89=2492", $true, $false, $false, $false, $false, $true, 1, $false, "44×93=4092", 2) | Out-Null
$d.Content.Find.Execute("57×91=5187", $true, $false, $false, $false, $false, $true, 1, $false, "23×23=529", 2) | Out-Null
$d.Content.Find.Execute("87×18=1566", $true, $false, $false, $false, $false, $true, 1, $false, "73×45=3285", 2) | Out-Null
